$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text storage (values like
# "41.439.36" or "0.0903" must stay literal strings, not be coerced to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '41.439.36'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '2.159.10'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '237.06'
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("D7").Value = '70.76'
$ws.Range("E7").Value = '  -1.86%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -3.96%  '
$ws.Range("D10").Value = '39.71'
$ws.Range("E10").Value = '  -4.47%  '
$ws.Range("D11").Value = '0.0903'
$ws.Range("E11").Value = '  -4.40%  '
$ws.Range("D12").Value = '54.28'
$ws.Range("E12").Value = '  -4.24%  '
$ws.Range("E13").Value = '  -3.31%  '
$ws.Range("E14").Value = '  -4.25%  '
$ws.Range("D15").Value = '2.482.03'
$ws.Range("E15").Value = '  -2.60%  '
$ws.Range("D16").Value = '14.27'
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").Value = '2.146.79'
$ws.Range("E17").Value = '  -3.62%  '
$ws.Range("E18").Value = '  -5.74%  '
$ws.Range("D19").Value = '41.278.45'
$ws.Range("E19").Value = '  -1.33%  '
$ws.Range("D20").Value = '0.0000101'
$ws.Range("E20").Value = '  -3.77%  '
$ws.Range("D21").Value = '69.33'
$ws.Range("E21").Value = '  -4.16%  '
$ws.Range("D22").Value = '5.75'
$ws.Range("E22").Value = '  -6.21%  '
$ws.Range("D23").Value = '9.76'
$ws.Range("E23").Value = '  -12.10%  '
$ws.Range("D24").Value = '225.11'
$ws.Range("E24").Value = '  -1.61%  '
$ws.Range("D25").Value = '1.97'
$ws.Range("E25").Value = '  -3.20%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("E27").Value = '  -6.55%  '
$ws.Range("D28").Value = '3.32'
$ws.Range("E28").Value = '  -8.40%  '
$ws.Range("E29").Value = '  -4.27%  '
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("D31").Value = '171.29'
$ws.Range("E31").Value = '  +2.52%  '
$ws.Range("D32").Value = '19.72'
$ws.Range("E32").Value = '  -3.09%  '
$ws.Range("D33").Value = '31.71'
$ws.Range("E33").Value = '  +6.52%  '
$ws.Range("D34").Value = '0.0764'
$ws.Range("E34").Value = '  -3.56%  '
$ws.Range("D35").Value = '5.09'
$ws.Range("E35").Value = '  -8.69%  '
$ws.Range("E36").Value = '  -3.31%  '
$ws.Range("D37").Value = '4.32'
$ws.Range("E37").Value = '  +3.09%  '
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("D40").Value = '12.15'
$ws.Range("E40").Value = '  -9.85%  '
$ws.Range("D41").Value = '2.05'
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("D42").Value = '5.31'
$ws.Range("E42").Value = '  -5.24%  '
$ws.Range("D43").Value = '58.22'
$ws.Range("E43").Value = '  -8.78%  '
$ws.Range("E44").Value = '  -4.10%  '
$ws.Range("D45").Value = '8.26'
$ws.Range("E45").Value = '  -4.55%  '
$ws.Range("E46").Value = '  -3.50%  '
$ws.Range("D47").Value = '97.60'
$ws.Range("E47").Value = '  -4.96%  '
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("E49").Value = '  -4.57%  '
$ws.Range("E50").Value = '  -7.37%  '
$ws.Range("E51").Value = '  -2.63%  '
